$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "67.457.50"
$ws.Cells.Item(2, 5).Value = "  -0.48%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.626.71"
$ws.Cells.Item(3, 5).Value = "  +0.70%  "
$ws.Cells.Item(4, 5).Value = "  -0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "592.61"
$ws.Cells.Item(5, 5).Value = "  -0.52%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "152.38"
$ws.Cells.Item(6, 5).Value = "  -1.96%  "
$ws.Cells.Item(7, 5).Value = "  +0.04%  "
$ws.Cells.Item(8, 5).Value = "  +0.88%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2.625.31"
$ws.Cells.Item(9, 5).Value = "  +0.63%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.123"
$ws.Cells.Item(10, 5).Value = "  -2.85%  "
$ws.Cells.Item(11, 5).Value = "  +0.18%  "
$ws.Cells.Item(12, 5).Value = "  -0.98%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.348"
$ws.Cells.Item(13, 5).Value = "  -1.70%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "27.70"
$ws.Cells.Item(14, 5).Value = "  +0.88%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.088.29"
$ws.Cells.Item(15, 5).Value = "  +0.11%  "
$ws.Cells.Item(16, 5).Value = "  -3.59%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "67.219.33"
$ws.Cells.Item(17, 5).Value = "  -0.71%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.623.14"
$ws.Cells.Item(18, 5).Value = "  +0.38%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "367.77"
$ws.Cells.Item(19, 5).Value = "  +0.65%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.09"
$ws.Cells.Item(20, 5).Value = "  -0.76%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "7.40"
$ws.Cells.Item(21, 5).Value = "  -2.83%  "
$ws.Cells.Item(23, 5).Value = "  -1.55%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.06"
$ws.Cells.Item(24, 5).Value = "  +3.04%  "
$ws.Cells.Item(25, 5).Value = "  +0.10%  "
$ws.Cells.Item(26, 5).Value = "  +2.66%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "67.63"
$ws.Cells.Item(27, 5).Value = "  +0.01%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.760.40"
$ws.Cells.Item(28, 5).Value = "  +0.67%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "585.82"
$ws.Cells.Item(29, 5).Value = "  +1.63%  "
$ws.Cells.Item(30, 5).Value = "  +0.13%  "
$ws.Cells.Item(31, 5).Value = "  -2.28%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.39"
$ws.Cells.Item(32, 5).Value = "  -2.44%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "7.72"
$ws.Cells.Item(33, 5).Value = "  -2.66%  "
$ws.Cells.Item(34, 5).Value = "  -2.62%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.00"
$ws.Cells.Item(35, 5).Value = "  -0.05%  "
$ws.Cells.Item(36, 5).Value = "  -7.01%  "
$ws.Cells.Item(37, 5).Value = "  -0.67%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "156.60"
$ws.Cells.Item(38, 5).Value = "  -1.14%  "
$ws.Cells.Item(39, 2).Value = "EthereumClassic"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "19.01"
$ws.Cells.Item(39, 5).Value = "  -1.71%  "
$ws.Cells.Item(40, 2).Value = "Stacks"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.88"
$ws.Cells.Item(40, 5).Value = "  +1.84%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.367"
$ws.Cells.Item(41, 5).Value = "  -0.53%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.25"
$ws.Cells.Item(42, 5).Value = "  -1.46%  "
$ws.Cells.Item(43, 5).Value = "  +2.04%  "
$ws.Cells.Item(44, 5).Value = "  +2.40%  "
$ws.Cells.Item(45, 5).Value = "  -0.01%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "155.78"
$ws.Cells.Item(46, 5).Value = "  +0.07%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0₆0296"
$ws.Cells.Item(47, 5).Value = "  +1.43%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.72"
$ws.Cells.Item(48, 5).Value = "  -0.19%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "21.93"
$ws.Cells.Item(49, 5).Value = "  +5.67%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.70"
$ws.Cells.Item(50, 5).Value = "  -0.79%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0786"
$ws.Cells.Item(51, 5).Value = "  +0.95%  "
